$wb = $excel.ActiveWorkbook

$navSheet = $wb.Worksheets.Item("NAV")
$navSheet.Range("A3").Value = "google_sheet"
$navSheet.Range("B3").Value = "1dc-SL4KNa9v89CE4lGxR1ZAdoyW1SbepHzKFf7I9__k"

$vehSheet = $wb.Worksheets.Item("5525")
$vehSheet.Range("A3").Value = "google_sheet"
$vehSheet.Range("B3").Value = "1ZDR9So-jv4lcPE9YTSno1Tde-xMdSHDDzReEIBmj55o"

$navSheet.Select()
$navSheet.Range("B3").Select()

$vehSheet.Activate()
$vehSheet.Range("B2").Select()

$navSheet.Activate()
